{"js": "// Fill in the user-story template fields:\n//  - story number \"#xx\" -> \"#02\"\n//  - collapse the long run of spaces before \"Titre\" down to a single space\n//  - title placeholder \"xx\" -> \"d\u00e9poser des vid\u00e9os de d\u00e9fis\"\n//  - priority \"1 \u00e0 5\" -> \"2\" (keeping the surrounding padding spaces)\n//  - remove the \"xx jours\" estimated-duration placeholder\n\nconst body = context.document.body;\n\n// 1) \"#xx\" -> \"#02\"\nconst storyNum = body.search(\"#xx\", { matchCase: true });\nstoryNum.load(\"items\");\nawait context.sync();\nif (storyNum.items.length > 0) {\n  storyNum.items[0].insertText(\"#02\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// 2) collapse the 26-space run that sits right before \"Titr\" down to one space\nconst spacesBeforeTitle = body.search(\" \".repeat(26) + \"Titr\", { matchCase: true });\nspacesBeforeTitle.load(\"items\");\nawait context.sync();\nif (spacesBeforeTitle.items.length > 0) {\n  spacesBeforeTitle.items[0].insertText(\" \" + \"Titr\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// 3) the title placeholder \"xx\" right after \"Titre : \" -> the real title\nconst titlePlaceholder = body.search(\"Titre\\u00a0: xx\", { matchCase: true });\ntitlePlaceholder.load(\"items\");\nawait context.sync();\nif (titlePlaceholder.items.length > 0) {\n  titlePlaceholder.items[0].insertText(\"Titre\\u00a0: d\u00e9poser des vid\u00e9os de d\u00e9fis\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// 4) \"1 \u00e0 5\" followed by 14 spaces -> \"2\" followed by 14 spaces (3+11 split, same total padding)\nconst priority = body.search(\"1 \u00e0 5\" + \" \".repeat(14), { matchCase: true });\npriority.load(\"items\");\nawait context.sync();\nif (priority.items.length > 0) {\n  priority.items[0].insertText(\"2\" + \" \".repeat(3) + \" \".repeat(11), Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// 5) drop the \"xx jours\" estimated-duration placeholder entirely\nconst duration = body.search(\"xx jours\", { matchCase: true });\nduration.load(\"items\");\nawait context.sync();\nif (duration.items.length > 0) {\n  duration.items[0].insertText(\"\", Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# Fill in the user-story template fields:\n#  - story number \"#xx\" -> \"#02\"\n#  - collapse the long run of spaces before \"Titre\" down to a single space\n#  - title placeholder \"xx\" -> \"d\u00e9poser des vid\u00e9os de d\u00e9fis\"\n#  - priority \"1 \u00e0 5\" -> \"2\" (keeping the surrounding padding spaces)\n#  - remove the \"xx jours\" estimated-duration placeholder\n\nfunction Repeat($s, $n) {\n    $out = \"\"\n    for ($i = 0; $i -lt $n; $i++) {\n        $out = $out + $s\n    }\n    return $out\n}\n\n$d = $word.ActiveDocument\n\nfunction DoReplace($findText, $replaceText) {\n    $r = $d.Content\n    $r.Find.ClearFormatting()\n    $r.Find.Replacement.ClearFormatting()\n    $res = $r.Find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 1)\n    return $res\n}\n\n# 1) \"#xx\" -> \"#02\"\nDoReplace \"#xx\" \"#02\" | Out-Null\n\n# 2) collapse the 26-space run right before \"Titr\" down to a single space\n$spaces26 = Repeat \" \" 26\n$find2 = $spaces26 + \"Titr\"\nDoReplace $find2 \" Titr\" | Out-Null\n\n# 3) the title placeholder \"xx\" right after \"Titre : \" -> the real title\n$nbsp = [char]0x00A0\n$find3 = \"Titre\" + $nbsp + \": xx\"\n$replace3 = \"Titre\" + $nbsp + \": d\u00e9poser des vid\u00e9os de d\u00e9fis\"\nDoReplace $find3 $replace3 | Out-Null\n\n# 4) \"1 \u00e0 5\" followed by 14 spaces -> \"2\" followed by 14 spaces (same total padding, split 3+11)\n$spaces14 = Repeat \" \" 14\n$spaces11 = Repeat \" \" 11\n$find4 = \"1 \u00e0 5\" + $spaces14\n$replace4 = \"2   \" + $spaces11\nDoReplace $find4 $replace4 | Out-Null\n\n# 5) drop the \"xx jours\" estimated-duration placeholder entirely\nDoReplace \"xx jours\" \"\" | Out-Null\n"}
